$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.532.71'
$ws.Range("E2").Value = '  +0.85%  '
$ws.Range("D3").Value = '2.378.74'
$ws.Range("E3").Value = '  +3.04%  '
$ws.Range("E4").Value = '  -0.06%  '
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '310.12'
$r.Style = "Normal"
$ws.Range("E5").Value = '  +0.02%  '
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '104.49'
$r.Style = "Normal"
$ws.Range("E6").Value = '  +3.59%  '
$ws.Range("E7").Value = '  -5.00%  '
$ws.Range("E8").Value = '  +0.01%  '
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '0.523'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -0.81%  '
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '35.96'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -0.09%  '
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '53.23'
$r.Style = "Normal"
$ws.Range("E11").Value = '  +1.77%  '
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.0812'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -0.71%  '
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.113'
$r.Style = "Normal"
$ws.Range("E13").Value = '  -0.50%  '
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '6.97'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -4.06%  '
$ws.Range("D15").Value = '2.748.45'
$ws.Range("E15").Value = '  +3.15%  '
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '15.62'
$r.Style = "Normal"
$ws.Range("E16").Value = '  +4.27%  '
$ws.Range("D17").Value = '2.382.23'
$ws.Range("E17").Value = '  +3.08%  '
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = '43.512.82'
$ws.Range("E19").Value = '  +0.96%  '
$ws.Range("E20").Value = '  +3.31%  '
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '11.90'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -5.16%  '
$ws.Range("E22").Value = '  -0.53%  '
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '68.35'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -0.04%  '
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '240.58'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -0.19%  '
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("E26").Value = '  -0.61%  '
$ws.Range("E27").Value = '  +0.16%  '
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '25.83'
$r.Style = "Normal"
$ws.Range("E28").Value = '  +4.58%  '
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = '3.86'
$r.Style = "Normal"
$ws.Range("E29").Value = '  -3.29%  '
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '36.69'
$r.Style = "Normal"
$ws.Range("E30").Value = '  -1.94%  '
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '9.51'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -1.46%  '
$ws.Range("E32").Value = '  -0.21%  '
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '160.96'
$r.Style = "Normal"
$ws.Range("E33").Value = '  -3.39%  '
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("E35").Value = '  -0.13%  '
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '18.28'
$r.Style = "Normal"
$ws.Range("E36").Value = '  +2.77%  '
$ws.Range("B37").Value = 'WEMIXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '2.53'
$r.Style = "Normal"
$ws.Range("E37").Value = '  +5.47%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '4.71'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +8.97%  '
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '3.10'
$r.Style = "Normal"
$ws.Range("E39").Value = '  -0.96%  '
$ws.Range("E40").Value = '  -0.53%  '
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '1.93'
$r.Style = "Normal"
$ws.Range("E41").Value = '  +5.26%  '
$ws.Range("E42").Value = '  -1.29%  '
$ws.Range("E43").Value = '  -2.07%  '
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '2.59'
$r.Style = "Normal"
$ws.Range("E44").Value = '  +13.11%  '
$ws.Range("D45").Value = '2.034.93'
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("E46").Value = '  +2.95%  '
$ws.Range("E47").Value = '  +0.46%  '
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '3.12'
$r.Style = "Normal"
$ws.Range("E48").Value = '  +3.24%  '
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = '10.57'
$r.Style = "Normal"
$ws.Range("E49").Value = '  +7.53%  '
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '57.95'
$r.Style = "Normal"
$ws.Range("E50").Value = '  +4.16%  '
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '2.95'
$r.Style = "Normal"
$ws.Range("E51").Value = '  -0.58%  '
